$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AD2").Value = 1000
$ws.Range("L12").Value = 1.2
$ws.Range("M12").Value = 4.33
$ws.Range("G15").Value = 2
$ws.Range("I15").Value = 3.25
$ws.Range("AC15").Value = 41
$ws.Range("AF15").Value = 17
$ws.Range("AI15").Value = 26
$ws.Range("G20").Value = 2.2
$ws.Range("H20").Value = 3.1
$ws.Range("I20").Value = 3.2
$ws.Range("L20").Value = 1.47
$ws.Range("M20").Value = 2.32
$ws.Range("N20").Value = 2.37
$ws.Range("O20").Value = 1.45
$ws.Range("P20").Value = 1.5
$ws.Range("Q20").Value = 2.25
$ws.Range("R20").Value = 2.05
$ws.Range("S20").Value = 1.6
$ws.Range("T20").Value = 5.8
$ws.Range("U20").Value = 9.25
$ws.Range("V20").Value = 9.75
$ws.Range("W20").Value = 21
$ws.Range("X20").Value = 22
$ws.Range("Y20").Value = 45
$ws.Range("Z20").Value = 6.8
$ws.Range("AA20").Value = 6.1
$ws.Range("AB20").Value = 19
$ws.Range("AC20").Value = 120
$ws.Range("AE20").Value = 7.4
$ws.Range("AF20").Value = 15
$ws.Range("AG20").Value = 12.5
$ws.Range("AH20").Value = 45
$ws.Range("AI20").Value = 35
$ws.Range("AJ20").Value = 55
$ws.Range("G24").Value = 1.88
$ws.Range("H24").Value = 3.05
$ws.Range("I24").Value = 4
$ws.Range("N24").Value = 2.12
$ws.Range("O24").Value = 1.57
$ws.Range("P24").Value = 1.42
$ws.Range("Q24").Value = 2.35
$ws.Range("T24").Value = 5.2
$ws.Range("U24").Value = 7
$ws.Range("V24").Value = 7.1
$ws.Range("W24").Value = 13
$ws.Range("X24").Value = 13.5
$ws.Range("Z24").Value = 7.4
$ws.Range("AA24").Value = 5.3
$ws.Range("AB24").Value = 13
$ws.Range("AC24").Value = 65
$ws.Range("AE24").Value = 8.25
$ws.Range("AF24").Value = 17
$ws.Range("AG24").Value = 11.25
$ws.Range("AH24").Value = 50
$ws.Range("AI24").Value = 32
$ws.Range("AJ24").Value = 40
$ws.Range("G26").Value = 3.15
$ws.Range("H26").Value = 3
$ws.Range("I26").Value = 2.2
$ws.Range("N26").Value = 2.15
$ws.Range("O26").Value = 1.55
$ws.Range("P26").Value = 1.44
$ws.Range("Q26").Value = 2.3
$ws.Range("T26").Value = 6.8
$ws.Range("U26").Value = 12.5
$ws.Range("V26").Value = 9.5
$ws.Range("W26").Value = 32
$ws.Range("X26").Value = 25
$ws.Range("Y26").Value = 32
$ws.Range("Z26").Value = 7.3
$ws.Range("AA26").Value = 5.1
$ws.Range("AB26").Value = 13
$ws.Range("AD26").Value = 500
$ws.Range("AE26").Value = 5.6
$ws.Range("AF26").Value = 8.25
$ws.Range("AG26").Value = 7.7
$ws.Range("AH26").Value = 17
$ws.Range("AI26").Value = 16
$ws.Range("AJ26").Value = 26
$ws.Range("G27").Value = 5.25
$ws.Range("H27").Value = 4
$ws.Range("I27").Value = 1.6
$ws.Range("R27").Value = 1.8
$ws.Range("S27").Value = 1.91
$ws.Range("U27").Value = 29
$ws.Range("V27").Value = 17
$ws.Range("Z27").Value = 12
$ws.Range("AE27").Value = 7.5
$ws.Range("AF27").Value = 8
$ws.Range("AH27").Value = 12
$ws.Range("AB29").Value = 14.5
$ws.Range("S30").Value = 1.82
$ws.Range("T30").Value = 6.3
$ws.Range("X30").Value = 16
$ws.Range("Y30").Value = 30
$ws.Range("AE30").Value = 10.25
$ws.Range("G32").Value = 1.7
$ws.Range("H32").Value = 3.55
$ws.Range("I32").Value = 4.85
$ws.Range("K32").Value = 7.2
$ws.Range("L32").Value = 1.3
$ws.Range("O32").Value = 1.83
$ws.Range("Q32").Value = 2.65
$ws.Range("S32").Value = 1.87
$ws.Range("U32").Value = 7.9
$ws.Range("Y32").Value = 26
$ws.Range("Z32").Value = 7.2
$ws.Range("AA32").Value = 6.9
$ws.Range("AF32").Value = 28
$ws.Range("AI32").Value = 50
$ws.Range("G33").Value = 2.42
$ws.Range("H33").Value = 2.87
$ws.Range("I33").Value = 3.15
$ws.Range("N33").Value = 2.42
$ws.Range("O33").Value = 1.5
$ws.Range("T33").Value = 6
$ws.Range("U33").Value = 10.5
$ws.Range("V33").Value = 10
$ws.Range("W33").Value = 26
$ws.Range("X33").Value = 25
$ws.Range("Y33").Value = 45
$ws.Range("AB33").Value = 17
$ws.Range("AE33").Value = 7.5
$ws.Range("AF33").Value = 15
$ws.Range("AG33").Value = 11.25
$ws.Range("AI33").Value = 32
$ws.Range("AJ33").Value = 45
$ws.Range("G38").Value = 11
$ws.Range("J38").Value = 23
$ws.Range("K38").Value = 1.02
$ws.Range("R38").Value = 1.83
$ws.Range("S38").Value = 1.83
$ws.Range("X38").Value = 67
$ws.Range("AE38").Value = 11
$ws.Range("AH38").Value = 8
$ws.Range("G39").Value = 1.33
$ws.Range("T39").Value = 9.5
$ws.Range("U39").Value = 7.5
$ws.Range("G40").Value = 1.83
$ws.Range("I40").Value = 3.75
$ws.Range("J40").Value = 1.02
$ws.Range("K40").Value = 12
$ws.Range("AI40").Value = 26
$ws.Range("G45").Value = 2.4
$ws.Range("I45").Value = 2.3
$ws.Range("L45").Value = 1.06
$ws.Range("M45").Value = 10
$ws.Range("N45").Value = 1.25
$ws.Range("O45").Value = 4
$ws.Range("R45").Value = 1.25
$ws.Range("S45").Value = 3.75
$ws.Range("V45").Value = 12
$ws.Range("AG45").Value = 11
